$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Action_recommendation_reference")

# Update PARTNER1's description in cell C10
$ws.Range("C10").Value = "NCC owns fee-simple land. Consider appropriate land sharing strategies."

# Move the active selection (matches the author's cursor position after editing)
$ws.Range("C11").Select()
